$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: column D (Price) and column E (Volume 1h) text values.
# Column D values that look like plain decimal numbers must be forced to
# stay as text (matching the original inline-string cell type) by applying
# a Text number format before assignment; otherwise Excel auto-converts
# them to numeric values and trailing/leading zeros would be lost.
$rowUpdates = @(
    [pscustomobject]@{ Row = 2; D = "30.602.42"; E = "  +0.58%  "; DForceText = $false }
    [pscustomobject]@{ Row = 3; D = "1.920.67"; E = "  -0.31%  "; DForceText = $false }
    [pscustomobject]@{ Row = 4; D = "0.9999"; E = "  -0.16%  "; DForceText = $true }
    [pscustomobject]@{ Row = 5; D = "247.76"; E = "  +3.20%  "; DForceText = $true }
    [pscustomobject]@{ Row = 6; D = "0.9999"; E = "  -0.24%  "; DForceText = $true }
    [pscustomobject]@{ Row = 7; D = "0.4744"; E = "  +0.23%  "; DForceText = $true }
    [pscustomobject]@{ Row = 8; D = $null; E = "  +1.52%  "; DForceText = $false }
    [pscustomobject]@{ Row = 9; D = "0.06842"; E = "  +4.07%  "; DForceText = $true }
    [pscustomobject]@{ Row = 10; D = "105.40"; E = "  +1.09%  "; DForceText = $true }
    [pscustomobject]@{ Row = 11; D = $null; E = "  -3.64%  "; DForceText = $false }
    [pscustomobject]@{ Row = 12; D = "1.920.88"; E = "  -0.28%  "; DForceText = $false }
    [pscustomobject]@{ Row = 13; D = "0.07698"; E = "  +1.51%  "; DForceText = $true }
    [pscustomobject]@{ Row = 14; D = "5.304"; E = "  +3.81%  "; DForceText = $true }
    [pscustomobject]@{ Row = 15; D = "0.6684"; E = "  +2.65%  "; DForceText = $true }
    [pscustomobject]@{ Row = 16; D = "291.93"; E = "  -1.21%  "; DForceText = $true }
    [pscustomobject]@{ Row = 17; D = "30.604.74"; E = "  +0.58%  "; DForceText = $false }
    [pscustomobject]@{ Row = 18; D = "0.000007601"; E = "  +1.45%  "; DForceText = $true }
    [pscustomobject]@{ Row = 19; D = "0.9996"; E = "  -0.33%  "; DForceText = $true }
    [pscustomobject]@{ Row = 20; D = $null; E = "  +0.05%  "; DForceText = $false }
    [pscustomobject]@{ Row = 21; D = "5.550"; E = "  +6.81%  "; DForceText = $true }
    [pscustomobject]@{ Row = 22; D = "2.169.97"; E = "  +0.09%  "; DForceText = $false }
    [pscustomobject]@{ Row = 23; D = $null; E = "  -0.12%  "; DForceText = $false }
    [pscustomobject]@{ Row = 24; D = "6.383"; E = "  +1.50%  "; DForceText = $true }
    [pscustomobject]@{ Row = 25; D = "9.411"; E = "  +1.72%  "; DForceText = $true }
    [pscustomobject]@{ Row = 26; D = "167.82"; E = "  +1.27%  "; DForceText = $true }
    [pscustomobject]@{ Row = 27; D = "21.14"; E = "  +8.39%  "; DForceText = $true }
    [pscustomobject]@{ Row = 28; D = "2.109"; E = "  +3.85%  "; DForceText = $true }
    [pscustomobject]@{ Row = 29; D = "0.1070"; E = "  -4.29%  "; DForceText = $true }
    [pscustomobject]@{ Row = 30; D = "1.395"; E = "  +2.78%  "; DForceText = $true }
    [pscustomobject]@{ Row = 31; D = "4.183"; E = "  +1.99%  "; DForceText = $true }
    [pscustomobject]@{ Row = 32; D = "4.069"; E = "  +3.72%  "; DForceText = $true }
    [pscustomobject]@{ Row = 33; D = "0.05031"; E = "  +0.42%  "; DForceText = $true }
    [pscustomobject]@{ Row = 34; D = "0.7396"; E = "  +0.34%  "; DForceText = $true }
    [pscustomobject]@{ Row = 35; D = "1.143"; E = "  +0.14%  "; DForceText = $true }
    [pscustomobject]@{ Row = 36; D = "0.02080"; E = "  +6.73%  "; DForceText = $true }
    [pscustomobject]@{ Row = 37; D = "2.745"; E = "  +0.78%  "; DForceText = $true }
    [pscustomobject]@{ Row = 38; D = "2.686"; E = "  -0.48%  "; DForceText = $true }
    [pscustomobject]@{ Row = 39; D = "2.057"; E = "  +1.96%  "; DForceText = $true }
    [pscustomobject]@{ Row = 40; D = "110.89"; E = "  +3.61%  "; DForceText = $true }
    [pscustomobject]@{ Row = 41; D = "0.8777"; E = "  +0.76%  "; DForceText = $true }
    [pscustomobject]@{ Row = 42; D = "0.4381"; E = "  +6.23%  "; DForceText = $true }
    [pscustomobject]@{ Row = 43; D = "5.880"; E = "  +0.80%  "; DForceText = $true }
    [pscustomobject]@{ Row = 44; D = "0.9998"; E = "  -0.30%  "; DForceText = $true }
    [pscustomobject]@{ Row = 45; D = "67.73"; E = "  -1.72%  "; DForceText = $true }
    [pscustomobject]@{ Row = 46; D = $null; E = "  +0.30%  "; DForceText = $false }
    [pscustomobject]@{ Row = 47; D = "9.382"; E = "  +2.24%  "; DForceText = $true }
    [pscustomobject]@{ Row = 48; D = "48.40"; E = "  +15.78%  "; DForceText = $true }
    [pscustomobject]@{ Row = 49; D = "0.1234"; E = "  +2.43%  "; DForceText = $true }
    [pscustomobject]@{ Row = 50; D = "34.84"; E = "  +0.71%  "; DForceText = $true }
    [pscustomobject]@{ Row = 51; D = "0.2480"; E = "  +10.90%  "; DForceText = $true }
)

foreach ($u in $rowUpdates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.DForceText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.Value = $u.E
}
